$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "NA" values under the duplicate_image_filename column (E) for data rows 2-21
for ($row = 2; $row -le 21; $row++) {
    $ws.Range("E$row").Value = "NA"
}
